# Generate Report for Handoff
# - Flip status from "In Translation" -> "Ready for handoff"
# - Bump the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" stamps
# - Widen the now-longer "Status"/locale-status columns to fit the new text

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-23 22:54:51"
$wsOverview.Range("E1").ColumnWidth = 16.3
$wsOverview.Range("F1").ColumnWidth = 16.3

# --- zh-cn sheet ------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-23 22:54:46"
$wsZhCn.Range("C1").ColumnWidth = 16.3

# --- de-de sheet ------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-23 22:54:51"
$wsDeDe.Range("C1").ColumnWidth = 16.3
